# Add the "CUI Flow Description" section to the end of the document,
# right before the final section properties (sectPr), mirroring the
# structure added in the target diff.

$d = $word.ActiveDocument

# Locate the very end of the document body content (after the last
# paragraph, "### If they delay evidence, your timeline slips.").
$lastParagraph = $d.Paragraphs.Last
$tail = $lastParagraph.Range
$insertionPoint = $d.Range($tail.End, $tail.End)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newContentXml = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>CUI Flow Description</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
  </w:pPr>
  <w:r>
    <w:t>Aegis Avionics receives Controlled Unclassified Information (CUI) from prime contractors via secure email, SharePoint portals, or other approved electronic transfer methods.</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Upon receipt, CUI is stored within designated SharePoint Online sites and accessed via authorized personnel using company-managed Windows endpoints. Users may synchronize CUI to local devices via OneDrive for Business, subject to endpoint security controls. </w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
  </w:pPr>
  <w:r>
    <w:t>CUI is processed by engineering and program staff using approved applications on managed endpoints. Internal collaboration involving CUI occurs through Microsoft Teams and SharePoint.</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
  </w:pPr>
  <w:r>
    <w:t>When CUI must be transmitted externally, it is shared through secure methods such as controlled SharePoint links or encrypted email, in accordance with company policy.</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Administrative access to systems handling CUI is restricted to authorized internal personnel and approved managed service provide (MSP) accounts, protected by multifactor authentication. </w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>CUI is not permitted on unmanaged or personal devices and nis not intentionally stored outside the defined CUI environment.</w:t>
  </w:r>
</w:p>
"@

$insertionPoint.InsertXML($newContentXml)
